$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

# Row 10 (1-based): columns N..S are 14..19
# N10 -> "70"   (new text, quote-prefixed so it's stored as text not a number)
# O10 -> "11"   (new text)
# P10 -> "12"   (reuses the existing shared string "12")
# Q10 -> "13"   (new text)
# R10 -> "Safe/Vault" (plain text, and its style is reset back to the default "Normal" style)
# S10 -> "211"  (new text)

$ws.Cells.Item(10, 14).Value2 = "'70"
$ws.Cells.Item(10, 15).Value2 = "'11"
$ws.Cells.Item(10, 16).Value2 = "'12"
$ws.Cells.Item(10, 17).Value2 = "'13"

$cellR10 = $ws.Cells.Item(10, 18)
$cellR10.Style = "Normal"
$cellR10.Value2 = "Safe/Vault"

$ws.Cells.Item(10, 19).Value2 = "'211"

# Update the sheet's active selection to match the new view state.
$ws.Activate()
$ws.Range("R21").Select()
